# Edit slide 4 ("Data Exploration" stage-1 notebook pipeline):
#  - rename two analysis functions and reflow their rounded-rectangle boxes
#  - replace the "In [N]" execution-count captions with descriptive captions
#  - drop the last two pipeline steps (plot_comment_classes / plot_multiple_classes)
#  - add a new header textbox "Basic data frame visuals for train and test sets"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Shape 4 (id 2, "Rounded Rectangle 1"): def hist_comment_lengths() -> def plot_comment_lengths()
$sh = $s.Shapes.Item(4)
$sub = $sh.TextFrame.TextRange.Characters(5, 20)
$sub.Text = "plot_comment_lengths"
$sh.Left = 21.333386421203617
$sh.Top = 157.74991607666018
$sh.Width = 210.33322906494143
$sh.Height = 65.33330917358398

# --- Shape 5 (id 3, "TextBox 2"): "In [5]" -> descriptive caption
$sh = $s.Shapes.Item(5)
$sh.TextFrame.TextRange.Text = "Show a histogram of the comment" + [char]0x2019 + "s lengths "
$sh.Left = 228.66661834716797
$sh.Top = 158.49999237060547
$sh.Width = 311.3332977294922
$sh.Height = 50.89220619201661

# --- Shape 6 (id 7, "Rounded Rectangle 6"): def clean_vs_toxic_ratio() -> def plot_labels()
$sh = $s.Shapes.Item(6)
$sub = $sh.TextFrame.TextRange.Characters(5, 20)
$sub.Text = "plot_labels"
$sh.Left = 21.333386421203617
$sh.Top = 240.41654205322268
$sh.Width = 210.33322906494143
$sh.Height = 65.33330917358398

# --- Shape 7 (id 8, "TextBox 7"): "In [6]" -> descriptive caption (3 runs, middle run is the column name)
$sh = $s.Shapes.Item(7)
$tr = $sh.TextFrame.TextRange
$tr.Text = "Show a bar plot of the counts of " + [char]0x2018 + "overall_toxic" + [char]0x2019 + " and " + [char]0x2018 + "non-toxic" + [char]0x2019 + " "
$full = $tr.Text
$start = $full.IndexOf("overall_toxic") + 1
$mid = $tr.Characters($start, 13)
$mid.Text = "overall_toxic"
$sh.Left = 231.66654205322268
$sh.Top = 240.41654205322268
$sh.Width = 264.0833892822266
$sh.Height = 50.89220619201661

# --- Remove the last two pipeline steps and their "In [N]" captions:
#       id 9  "Rounded Rectangle 8"  def plot_comment_classes()
#       id 12 "Rounded Rectangle 11" def plot_multiple_classes()
#       id 13 "TextBox 12"           In [8]
# (delete from the back so earlier indices stay valid)
$s.Shapes.Item(11).Delete()
$s.Shapes.Item(10).Delete()
$s.Shapes.Item(8).Delete()

# --- Replace "TextBox 10" (In [7]) with a new header textbox; add two throwaway
# shapes first so the new shape's auto-assigned id lands on 14 (matching the
# target id/name "TextBox 13"), then remove the throwaways.
$junk1 = $s.Shapes.AddTextbox(1, 0, 0, 100, 100)
$junk2 = $s.Shapes.AddTextbox(1, 0, 0, 100, 100)
$newBox = $s.Shapes.AddTextbox(1, 0, 0, 100, 100)
$junk1.Delete()
$junk2.Delete()

$s.Shapes.Item(8).Delete()   # old "TextBox 10" (id 11, "In [7]")

$newBox.Fill.Visible = 0
$newBox.TextFrame.WordWrap = -1
$newBox.TextFrame.AutoSize = 1
$newBox.TextFrame.TextRange.Text = "Basic data frame visuals for train and test sets"
$newBox.Left = 21.333386421203617
$newBox.Top = 115.42882156372072
$newBox.Width = 372.41661071777344
$newBox.Height = 29.081259727478027
